$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 273, shifting existing rows 273:356 down to 274:357
$ws.Rows(273).Insert()

# Populate the newly inserted row 273 with the new data record
$ws.Range("A273").Value = 4
$ws.Range("B273").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C273").Value = "Los Lagos"
$ws.Range("D273").Value = 44876
$ws.Range("E273").Value = 10
$ws.Range("F273").Value = 100112037
$ws.Range("G273").Value = "Cebollín"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 160
$ws.Range("K273").Value = 6000
$ws.Range("L273").Value = 6000
$ws.Range("M273").Value = 6000
$ws.Range("N273").Value = "$/paquete 36 unidades"
$ws.Range("O273").Value = "Región Metropolitana"
$ws.Range("P273").Value = 167
$ws.Range("Q273").Value = 36
$ws.Range("R273").Value = "Hortaliza"
